# Clear the three quarterly estimate rows (2024-09, 2024-12, 2025-03 periods),
# leaving only row 24 (2024-06 period) populated. The content (chen/aa/date
# range/unit price/qty/subtotal) is removed and the cells revert to the plain
# row style (same as B25:B27 / the blank rows 28-29), so the dependent totals
# (H31 subtotal, H32 tax, H33 total, C21 grand total) recalculate down from
# 4 populated rows to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the values/formulas from C25:H27.
$ws.Range("C25:H27").ClearContents()

# Re-stamp those cells with the same plain style used by column B in those
# rows (style index 24), matching how the now-empty rows look.
$ws.Range("B25").Copy()
$ws.Range("C25:H27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
